# Clean up data columns to conform to specs:
#  - experimentDesign (column D): "Environmental Perturbation" -> "Environmental_Perturbation"
#  - strain (column F): "KN99 alpha" -> "KN99_alpha"
# Also update the sheet's active selection from B2:B27 to F2:F27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 27 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    if ($dCell.Value2 -eq "Environmental Perturbation") {
        $dCell.Value2 = "Environmental_Perturbation"
    }

    $fCell = $ws.Cells.Item($r, 6)
    if ($fCell.Value2 -eq "KN99 alpha") {
        $fCell.Value2 = "KN99_alpha"
    }
}

$ws.Range("F2:F27").Select()
